$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("default")

# The password for user in row 11 (A11 = "vertexlc") was changed
# from the weak/common password "letmein1" to a stronger one "dm8k26h2".
$ws.Range("B11").Value = "dm8k26h2"
